# ===================================================================
# Edit script for ZBabcock_CritiqueKenAlbala.docx
# Applies the changes described by the authoritative XML diff.
# ===================================================================

$d = $word.ActiveDocument

$LDQ = [char]0x201C   # “
$RDQ = [char]0x201D   # ”
$RSQ = [char]0x2019   # ’

function Replace-InParagraph {
    param(
        [int]$ParaIndex,
        [string]$OldText,
        [string]$NewText
    )
    $p = $d.Paragraphs($ParaIndex)
    $searchRange = $d.Range($p.Range.Start, $p.Range.End)
    $f = $searchRange.Find
    $f.ClearFormatting()
    $found = $f.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)
    return $found
}

# -------------------------------------------------------------
# 1) Paragraph 9 (index): replace strike-through sentence at the
#    end with the new "counterpoint" text, dropping the strike
#    formatting entirely.
# -------------------------------------------------------------
$p9 = $d.Paragraphs(9)
$p9Start = $p9.Range.Start
$p9EndAll = $p9.Range.End
$searchRange = $d.Range($p9Start, $p9EndAll)
$f = $searchRange.Find
$f.ClearFormatting()
$found = $f.Execute(", all without meandering.")
$startReplace = $searchRange.Start
$delRange = $d.Range($startReplace, $p9EndAll - 1)
$delRange.Text = ""

$ins = $d.Range($startReplace, $startReplace)
$ins.InsertAfter(", all without meandering.")
$ins.Collapse(0)
$ins.InsertAfter(" ")
$ins.Collapse(0)
$ins.InsertAfter(" ")
$ins.Collapse(0)
$ins.InsertAfter("Albala")
$ins.Collapse(0)
$ins.InsertAfter(" does, however, raise one counterpoint to his")
$ins.Collapse(0)
$ins.InsertAfter(" own")
$ins.Collapse(0)
$ins.InsertAfter(" that seems to be quite obscure: ")
$ins.Collapse(0)
$ins.InsertAfter("that " + $LDQ + "people are genetically adapted to certain regions and its local food over time," + $RDQ + " just like physical features through ")
$ins.Collapse(0)
$ins.InsertAfter("evolution.")
$ins.Collapse(0)
$ins.InsertAfter(" While this point lends itself moderately well to his following point, of how ingredients have been moving around for thousands of years, it perhaps would" + $RSQ + "ve been more beneficial to have brought up a more widely-held point of view.")
$ins.Collapse(0)
$ins.InsertAfter(" ")

Write-Output "Step 1 done"

# -------------------------------------------------------------
# 2) Paragraph 10: "Albala, however, does seem to have missed an
#    important aspect that may be a strong foothold for
#    authenticity though: the aspect of health." ->
#    "Albala does seem to have missed a more universal, stable
#    concern with food authenticity: the aspect of health."
# -------------------------------------------------------------
$ok = Replace-InParagraph 10 "Albala, however, does seem to have missed an important aspect that may be a strong foothold for authenticity though: the aspect of health." "Albala does seem to have missed a more universal, stable concern with food authenticity: the aspect of health."
Write-Output ("Step 2: " + $ok)

# -------------------------------------------------------------
# 3) Paragraph 10: add the missing sub-question to the Lupo
#    article title quoted in-line.
# -------------------------------------------------------------
$old3 = $LDQ + "Food Authenticity:" + $RDQ + ", Lisa Lupo"
$new3 = $LDQ + "Food Authenticity: What is spurring this anti-fraud movement?" + $RDQ + ", Lisa Lupo"
$ok = Replace-InParagraph 10 $old3 $new3
Write-Output ("Step 3: " + $ok)

# -------------------------------------------------------------
# 4) Paragraph 10: insert "[by humans]" into the melamine quote.
# -------------------------------------------------------------
$old4 = $LDQ + "ingested in large doses may cause stones and illness"
$new4 = $LDQ + "ingested in large doses [by humans] may cause stones and illness"
$ok = Replace-InParagraph 10 $old4 $new4
Write-Output ("Step 4: " + $ok)

# -------------------------------------------------------------
# 5) Paragraph 10: extend the final sentence.
# -------------------------------------------------------------
$p10 = $d.Paragraphs(10)
$searchRange = $d.Range($p10.Range.Start, $p10.Range.End)
$f = $searchRange.Find
$f.ClearFormatting()
$found = $f.Execute(", due to the significant health risks that may follow. ")
Write-Output ("Step 5 find: " + $found)
$matchStart = $searchRange.Start
$matchEnd = $searchRange.End
$d.Range($matchStart, $matchEnd).Text = ""
$ins5 = $d.Range($matchStart, $matchStart)
$ins5.InsertAfter(", due to the significant health risks that may follow")
$ins5.Collapse(0)
$ins5.InsertAfter(", which is something Albala should" + $RSQ + "ve addressed.")
Write-Output "Step 5 done"

# -------------------------------------------------------------
# 6) Paragraph 11 ("With this abstract concept ...") gets a
#    wholly new closing-argument paragraph, the stray "_GoBack"
#    bookmark there is removed, and the following (page-break-
#    only) paragraph is folded into this one.
# -------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$p11 = $d.Paragraphs(11)
$tabEnd = $p11.Range.Start + 1
$textEnd = $p11.Range.End - 1
$d.Range($tabEnd, $textEnd).Text = ""

$ins6 = $d.Range($tabEnd, $tabEnd)
$ins6.InsertAfter("Authenticity does have a place in our society")
$ins6.Collapse(0)
$ins6.InsertAfter(". T")
$ins6.Collapse(0)
$ins6.InsertAfter("he")
$ins6.Collapse(0)
$ins6.InsertAfter(" problem that has arisen is that it" + $RSQ + "s had too much weight put on it")
$ins6.Collapse(0)
$ins6.InsertAfter(", in faulty ways")
$ins6.Collapse(0)
$ins6.InsertAfter(". Albala" + $RSQ + "s address of this")
$ins6.Collapse(0)
$ins6.InsertAfter(" problem")
$ins6.Collapse(0)
$ins6.InsertAfter(" ")
$ins6.Collapse(0)
$ins6.InsertAfter("shows his thorough knowledge of the subject")
$ins6.Collapse(0)
$ins6.InsertAfter(", ")
$ins6.Collapse(0)
$ins6.InsertAfter("but he took his conclusion a step too far.")
$ins6.Collapse(0)
$ins6.InsertAfter(" ")
$ins6.Collapse(0)
$ins6.InsertAfter("Perhaps if his proposal of doing away with the concept of authenticity was executed, " + $LDQ + "food authenticity" + $RDQ + " as Karen ")
$ins6.Collapse(0)
$ins6.InsertAfter("Everstein")
$ins6.Collapse(0)
$ins6.InsertAfter(" define")
$ins6.Collapse(0)
$ins6.InsertAfter("s")
$ins6.Collapse(0)
$ins6.InsertAfter(" it would just be referred to by another name. If this was guaranteed to happen, Albala" + $RSQ + "s article would")
$ins6.Collapse(0)
$ins6.InsertAfter(" ")
$ins6.Collapse(0)
$ins6.InsertAfter("be a definite must-read that many people would benefit greatly from")
$ins6.Collapse(0)
$ins6.InsertAfter(".")
$ins6.Collapse(0)
$ins6.InsertAfter(" ")
$ins6.Collapse(0)
$ins6.InsertAfter("The problem is that")
$ins6.Collapse(0)
$ins6.InsertAfter(" it" + $RSQ + "s uncertain whether ")
$ins6.Collapse(0)
$ins6.InsertAfter("Everstein" + $RSQ + "s")
$ins6.Collapse(0)
$ins6.InsertAfter(" definition would be caught up in a discreditation of food authenticity")
$ins6.Collapse(0)
$ins6.InsertAfter(". ")
$ins6.Collapse(0)
$ins6.InsertAfter("Albala" + $RSQ + "s piece is still an article ")
$ins6.Collapse(0)
$ins6.InsertAfter("that many would benefit greatly from")
$ins6.Collapse(0)
$ins6.InsertAfter(", but there" + $RSQ + "s an underlying problem with their knowledge")
$ins6.Collapse(0)
$ins6.InsertAfter(":")
$ins6.Collapse(0)
$ins6.InsertAfter(" ")
$ins6.Collapse(0)
$ins6.InsertAfter("many")
$ins6.Collapse(0)
$ins6.InsertAfter(" ")
$ins6.Collapse(0)
$ins6.InsertAfter("people")
$ins6.Collapse(0)
$ins6.InsertAfter(" could")
$ins6.Collapse(0)
$ins6.InsertAfter(" ")
$ins6.Collapse(0)
$ins6.InsertAfter("develop ")
$ins6.Collapse(0)
$ins6.InsertAfter("a harmful dichotomous viewpoint that could make food regulations viewed as mere suggestions")
$ins6.Collapse(0)
$ins6.InsertAfter(". This is ")
$ins6.Collapse(0)
$ins6.InsertAfter("something that is already a bad enough problem in the United States,")
$ins6.Collapse(0)
$ins6.InsertAfter(" and")
$ins6.Collapse(0)
$ins6.InsertAfter(" one that needs to be addressed")
$ins6.Collapse(0)
$ins6.InsertAfter(" as well. Albala has a point, but people need to be careful to not make matters worse.")

# fold the page-break-only paragraph into this one
$p11 = $d.Paragraphs(11)
$endOfP11 = $p11.Range.End
$d.Range($endOfP11 - 1, $endOfP11).Delete()

Write-Output "Step 6 done"
Write-Output $d.Paragraphs(11).Range.Text
